# "TP Location for debug"
# Adds the new Test-Point / net-name lookup table (rows 32-54, columns C & E)
# to the "Airzai Interconnects" sheet, and leaves the final selection on E46,
# matching the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered so that brand-new shared strings are introduced in the same
# sequence the author typed them in (TP6..TP22 down column C first, then
# the column E net names, then TP23, then the remaining INT net names).
$newCells = [ordered]@{
  "C32" = "TP6"
  "C33" = "TP7"
  "C34" = "TP8"
  "C35" = "TP9"
  "C36" = "TP10"
  "C37" = "TP11"
  "C38" = "TP12"
  "C39" = "TP13"
  "C40" = "TP14"
  "C41" = "TP15"
  "C42" = "TP16"
  "C43" = "TP17"
  "C44" = "TP18"
  "C45" = "TP19"
  "C46" = "TP20"
  "C47" = "TP21"
  "C48" = "TP22"
  "E33" = "VCC_5V_USB"
  "E32" = "VCC_3V_Reg"
  "E36" = "VCC_3V_J14"
  "E37" = "VCC_5V_J14"
  "E38" = "EN_Line"
  "E39" = "LED_STR_Data"
  "E40" = "LED_STR_CLK"
  "C49" = "TP23"
  "E42" = "TOF_INT"
  "E43" = "NFC_INT"
  "E44" = "ACCL_INT"
  "E34" = "I2C_SCL"
  "E35" = "I2C_SDA"
  "E45" = "PWM1"
  "C50" = "TP18"
  "C51" = "TP19"
  "C52" = "TP20"
  "C53" = "TP21"
  "C54" = "TP22"
}

foreach ($addr in $newCells.Keys) {
    $ws.Range($addr).Value2 = $newCells[$addr]
}

# Final cursor/selection position left by the author when the workbook was saved.
$ws.Range("E46").Select()
